$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1943.8889
$ws.Range("I40").Value = 1436.875
$ws.Range("K40").Value = 1436.875
$ws.Range("M40").Value = -1261.875
$ws.Range("H70").Value = 3250.25
$ws.Range("I70").Value = 3250.25
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9750.75
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9480.75
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 3250.25
$ws.Range("I73").Value = 3250.25
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9750.75
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8814.75
$ws.Range("N73").ClearContents()
$ws.Range("H99").Value = 732.7
$ws.Range("I99").Value = 618.5714
$ws.Range("K99").Value = 1855.7142
$ws.Range("M99").Value = -357.7142000000001
$ws.Range("H113").Value = 6645.3335
$ws.Range("I113").Value = 5713.857
$ws.Range("J113").Value = 7460.375
$ws.Range("K113").Value = 5713.857
$ws.Range("L113").Value = 7460.375
$ws.Range("M113").Value = -2459.857
$ws.Range("N113").Value = -13968.375
$ws.Range("H138").Value = 3569.0476
$ws.Range("J138").Value = 3777.7778
$ws.Range("L138").Value = 11333.3334
$ws.Range("N138").Value = -21613.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H45").Value = 1664
$ws.Range("I45").Value = 1664
$ws.Range("K45").Value = 1664
$ws.Range("M45").Value = -1287
$ws.Range("H53").Value = 10000
$ws.Range("J53").Value = 10000
$ws.Range("L53").Value = 10000
$ws.Range("N53").Value = -11364
$ws.Range("H61").Value = 2703.6875
$ws.Range("I61").Value = 2616.5557
$ws.Range("J61").Value = 2815.7144
$ws.Range("K61").Value = 2616.5557
$ws.Range("L61").Value = 2815.7144
$ws.Range("M61").Value = -2404.5557
$ws.Range("N61").Value = -3239.7144
$ws.Range("H63").Value = 1602.5555
$ws.Range("J63").Value = 4000
$ws.Range("L63").Value = 4000
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 1602.5555
$ws.Range("J66").Value = 4000
$ws.Range("L66").Value = 20000
$ws.Range("N66").Value = -26864
$ws.Range("H74").Value = 2664.3333
$ws.Range("I74").Value = 2497.5
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2497.5
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1623.5
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 2664.3333
$ws.Range("I77").Value = 2497.5
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 12487.5
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -8119.5
$ws.Range("N77").Value = -33736
$ws.Range("H122").Value = 2723.6667
$ws.Range("I122").Value = 2723.6667
$ws.Range("K122").Value = 8171.000100000001
$ws.Range("M122").Value = -5721.000100000001
$ws.Range("H136").Value = 2703.6875
$ws.Range("I136").Value = 2616.5557
$ws.Range("J136").Value = 2815.7144
$ws.Range("K136").Value = 7849.6671
$ws.Range("L136").Value = 8447.143199999999
$ws.Range("M136").Value = -5299.6671
$ws.Range("N136").Value = -13547.1432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5911.4346
$ws.Range("I134").Value = 6295.778
$ws.Range("J134").Value = 4527.8
$ws.Range("K134").Value = 18887.334
$ws.Range("L134").Value = 13583.4
$ws.Range("M134").Value = -16352.334
$ws.Range("N134").Value = -18653.4
$ws.Range("H135").Value = 45997.6
$ws.Range("J135").Value = 45997.6
$ws.Range("L135").Value = 45997.6
$ws.Range("N135").Value = -56137.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1100
$ws.Range("H31").Value = 1966.375
$ws.Range("I31").Value = 1790.1428
$ws.Range("J31").Value = 3200
$ws.Range("K31").Value = 1790.1428
$ws.Range("L31").Value = 3200
$ws.Range("M31").Value = -1495.1428
$ws.Range("N31").Value = -3790
$ws.Range("H34").Value = 1966.375
$ws.Range("I34").Value = 1790.1428
$ws.Range("J34").Value = 3200
$ws.Range("K34").Value = 1790.1428
$ws.Range("L34").Value = 3200
$ws.Range("M34").Value = -1588.1428
$ws.Range("N34").Value = -3604
$ws.Range("H62").Value = 1700.6666
$ws.Range("I62").Value = 1601
$ws.Range("J62").Value = 1900
$ws.Range("K62").Value = 1601
$ws.Range("L62").Value = 1900
$ws.Range("M62").Value = -977
$ws.Range("N62").Value = -3148
$ws.Range("H65").Value = 1700.6666
$ws.Range("I65").Value = 1601
$ws.Range("J65").Value = 1900
$ws.Range("K65").Value = 8005
$ws.Range("L65").Value = 9500
$ws.Range("M65").Value = -4885
$ws.Range("N65").Value = -15740
$ws.Range("H141").Value = 40000
$ws.Range("I141").Value = 40000
$ws.Range("K141").Value = 40000
$ws.Range("M141").Value = -34820

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1475.2963
$ws.Range("I4").Value = 1391.2106
$ws.Range("K4").Value = 4173.6318
$ws.Range("M4").Value = -4061.6318
$ws.Range("H5").Value = 1277.4
$ws.Range("J5").Value = 900
$ws.Range("L5").Value = 2700
$ws.Range("N5").Value = -2924
$ws.Range("H48").Value = 400
$ws.Range("I48").Value = 400
$ws.Range("K48").Value = 1200
$ws.Range("M48").Value = -950
$ws.Range("H92").Value = 994
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 1123.3334
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 3370.0002
$ws.Range("M92").Value = -1152
$ws.Range("N92").Value = -5866.0002
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H135").Value = 1277.4
$ws.Range("J135").Value = 900
$ws.Range("L135").Value = 8100
$ws.Range("N135").Value = -13170

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2499
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 17252.092
$ws.Range("J16").Value = 77496
$ws.Range("L16").Value = 77496
$ws.Range("N16").Value = -77836
$ws.Range("H55").Value = 332.25
$ws.Range("I55").Value = 303.83334
$ws.Range("K55").Value = 303.83334
$ws.Range("M55").Value = -130.83334
$ws.Range("H68").Value = 73400.2
$ws.Range("J68").Value = 119999.664
$ws.Range("L68").Value = 119999.664
$ws.Range("N68").Value = -121497.664
$ws.Range("H71").Value = 73400.2
$ws.Range("J71").Value = 119999.664
$ws.Range("L71").Value = 599998.3200000001
$ws.Range("N71").Value = -607486.3200000001
$ws.Range("H93").Value = 2934.3333
$ws.Range("I93").Value = 2934.3333
$ws.Range("K93").Value = 2934.3333
$ws.Range("M93").Value = -1686.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 7500
$ws.Range("J48").Value = 7500
$ws.Range("L48").Value = 7500
$ws.Range("N48").Value = -8638
$ws.Range("H107").Value = 1167.5333
$ws.Range("I107").Value = 837.5454999999999
$ws.Range("J107").Value = 2075
$ws.Range("K107").Value = 2512.6365
$ws.Range("L107").Value = 6225
$ws.Range("M107").Value = -592.6364999999996
$ws.Range("N107").Value = -10065
